# Update cryptos list with latest price/volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to stay plain text even when the string looks numeric
    # (mirrors how the source data was stored as inline text), then drop the
    # temporary "@" number format so the cell style is left untouched.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.472.92"
$ws.Range("E2").Value = "  +0.21%  "

Set-TextValue $ws.Range("D3") "1.820.12"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("E4").Value = "  -0.09%  "

Set-TextValue $ws.Range("D5") "315.05"
$ws.Range("E5").Value = "  -0.41%  "

$ws.Range("E6").Value = "  -0.09%  "

Set-TextValue $ws.Range("D7") "0.5095"
$ws.Range("E7").Value = "  -4.80%  "

Set-TextValue $ws.Range("D8") "0.3952"
$ws.Range("E8").Value = "  -1.73%  "

Set-TextValue $ws.Range("D9") "0.08273"
$ws.Range("E9").Value = "  +7.96%  "

$ws.Range("E10").Value = "  -0.07%  "

Set-TextValue $ws.Range("D11") "41.57"
$ws.Range("E11").Value = "  -0.75%  "

Set-TextValue $ws.Range("D13") "21.10"
$ws.Range("E13").Value = "  +0.87%  "

$ws.Range("E14").Value = "  -0.03%  "

Set-TextValue $ws.Range("D15") "7.530"
$ws.Range("E15").Value = "  -1.54%  "

Set-TextValue $ws.Range("D16") "1.817.13"
$ws.Range("E16").Value = "  -0.44%  "

Set-TextValue $ws.Range("D17") "0.00001146"
$ws.Range("E17").Value = "  +6.36%  "

Set-TextValue $ws.Range("D18") "92.57"
$ws.Range("E18").Value = "  +3.29%  "

Set-TextValue $ws.Range("D19") "0.06657"
$ws.Range("E19").Value = "  +1.09%  "

Set-TextValue $ws.Range("D20") "17.77"
$ws.Range("E20").Value = "  +0.66%  "

$ws.Range("E21").Value = "  -0.05%  "

Set-TextValue $ws.Range("D22") "6.117"
$ws.Range("E22").Value = "  +0.74%  "

Set-TextValue $ws.Range("D23") "28.477.79"
$ws.Range("E23").Value = "  +0.19%  "

Set-TextValue $ws.Range("D24") "11.45"
$ws.Range("E24").Value = "  +3.26%  "

Set-TextValue $ws.Range("D25") "2.269"
$ws.Range("E25").Value = "  +2.08%  "

Set-TextValue $ws.Range("D26") "21.31"
$ws.Range("E26").Value = "  +2.95%  "

Set-TextValue $ws.Range("D27") "155.93"
$ws.Range("E27").Value = "  -1.06%  "

Set-TextValue $ws.Range("D28") "2.025.47"
$ws.Range("E28").Value = "  -0.62%  "

Set-TextValue $ws.Range("D29") "2.409"
$ws.Range("E29").Value = "  -1.92%  "

Set-TextValue $ws.Range("D30") "125.82"
$ws.Range("E30").Value = "  +1.61%  "

Set-TextValue $ws.Range("D31") "1.108"
$ws.Range("E31").Value = "  -1.16%  "

Set-TextValue $ws.Range("D32") "0.1092"
$ws.Range("E32").Value = "  -2.03%  "

Set-TextValue $ws.Range("D33") "5.788"
$ws.Range("E33").Value = "  +2.02%  "

Set-TextValue $ws.Range("D34") "3.650"
$ws.Range("E34").Value = "  +0.17%  "

Set-TextValue $ws.Range("D35") "0.07071"
$ws.Range("E35").Value = "  -4.13%  "

$ws.Range("E36").Value = "  -0.75%  "

Set-TextValue $ws.Range("D37") "0.02339"
$ws.Range("E37").Value = "  -0.09%  "

Set-TextValue $ws.Range("D38") "5.231"
$ws.Range("E38").Value = "  +0.07%  "

Set-TextValue $ws.Range("D39") "8.853"
$ws.Range("E39").Value = "  +0.02%  "

Set-TextValue $ws.Range("D40") "0.6292"
$ws.Range("E40").Value = "  +0.40%  "

Set-TextValue $ws.Range("D41") "11.30"

Set-TextValue $ws.Range("D42") "1.180"
$ws.Range("E42").Value = "  -0.01%  "

Set-TextValue $ws.Range("D43") "1.0000"

Set-TextValue $ws.Range("D44") "1.400"
$ws.Range("E44").Value = "  +0.50%  "

Set-TextValue $ws.Range("D45") "13.47"
$ws.Range("E45").Value = "  -0.06%  "

Set-TextValue $ws.Range("D46") "0.5914"
$ws.Range("E46").Value = "  +1.24%  "

Set-TextValue $ws.Range("D47") "3.730"
$ws.Range("E47").Value = "  +0.90%  "

Set-TextValue $ws.Range("D48") "125.23"
$ws.Range("E48").Value = "  +0.42%  "

Set-TextValue $ws.Range("D49") "1.984"
$ws.Range("E49").Value = "  -1.01%  "

$ws.Range("E50").Value = "  -1.49%  "

Set-TextValue $ws.Range("D51") "0.06890"
$ws.Range("E51").Value = "  +0.07%  "
